# BVT checklist update for Circular Gauge visual (22 visuals BVT refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "BVTs" sheet

# --- Append two new BVT rows covering Actual/Target value edge cases -------
$ws.Range("B21").Value = "Actual Value"
$ws.Range("B22").Value = "Target Value"

$ws.Range("C21").Value = "Actual value accept 0 values but negative values"
$ws.Range("C22").Value = "Target value accept 0 values but negative values"

$ws.Range("D21").Value = "1.Acual Value =0,Target Value=0                                            2.Actual Value=0,Target Value> Actual Value"
$ws.Range("D22").Value = "1.Target Value=0,Target Value< Actual Value    "
$ws.Range("D21:D22").WrapText = $true

$ws.Range("E21").Value = "1.Visual shows 0%                                                                                                                         2.Visual shows 0%            "
$ws.Range("E22").Value = "1.Visual shows 100%            "
$ws.Range("E21").WrapText = $true

$ws.Range("A21:E21").RowHeight = 74.25

# --- Update the "Drag columns" step descriptions in rows 2 & 3 --------------
$ws.Range("D2").Value = "1. Drag 'Last Year Sales' column in 'Actual Value' field`n2. Drag 'Target' column in 'Target Value' field and summarize it to sum"
$ws.Range("D3").Value = "1. Drag 'Target' column in 'Tooltip' field`n2. Hover over the circle in the visual to see the tooltips"

# Row 2 grew from a 2-line to a 4-line instruction, so its row height doubles.
$ws.Range("A2:E2").RowHeight = 60

# --- Page setup (sheet now prints as A4 portrait) ---------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Restore the view: scrolled to the top, selection on D4 -----------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("D4").Select()
